$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2095808383233533
$ws.Range("C2").Value = 0.5179640718562875
$ws.Range("J2").Value = 0.02095808383233533
$ws.Range("O2").Value = 0.002994011976047904
$ws.Range("P2").Value = 0.1467065868263473
$ws.Range("S2").Value = 0.1017964071856287

$ws.Range("B3").Value = 0.0111731843575419
$ws.Range("C3").Value = 0.0335195530726257
$ws.Range("J3").Value = 0.03910614525139665
$ws.Range("P3").Value = 0.6983240223463687
$ws.Range("S3").Value = 0.217877094972067

$ws.Range("J4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.7872340425531915
$ws.Range("S4").Value = 0.1914893617021277

$ws.Range("B6").Value = 0.0593607305936073
$ws.Range("D6").Value = 0.0091324200913242
$ws.Range("F6").Value = 0.0502283105022831
$ws.Range("J6").Value = 0.273972602739726
$ws.Range("O6").Value = 0.0273972602739726
$ws.Range("Q6").Value = 0.1598173515981735
$ws.Range("R6").Value = 0.0730593607305936
$ws.Range("S6").Value = 0.3470319634703196

$ws.Range("B7").Value = 0.1131221719457014
$ws.Range("D7").Value = 0.02714932126696833
$ws.Range("E7").Value = 0.004524886877828055
$ws.Range("F7").Value = 0.04072398190045249
$ws.Range("J7").Value = 0.1312217194570136
$ws.Range("O7").Value = 0.01357466063348416
$ws.Range("Q7").Value = 0.1447963800904978
$ws.Range("R7").Value = 0.08597285067873303
$ws.Range("S7").Value = 0.4389140271493213

$ws.Range("B8").Value = 0.0958904109589041
$ws.Range("D8").Value = 0.0182648401826484
$ws.Range("F8").Value = 0.0730593607305936
$ws.Range("J8").Value = 0.1415525114155251
$ws.Range("O8").Value = 0.01141552511415525
$ws.Range("Q8").Value = 0.136986301369863
$ws.Range("R8").Value = 0.06164383561643835
$ws.Range("S8").Value = 0.4611872146118721

$ws.Range("B9").Value = 0.105726872246696
$ws.Range("D9").Value = 0.01762114537444934
$ws.Range("E9").Value = 0.004405286343612335
$ws.Range("F9").Value = 0.05726872246696035
$ws.Range("J9").Value = 0.105726872246696
$ws.Range("O9").Value = 0.004405286343612335
$ws.Range("Q9").Value = 0.1497797356828194
$ws.Range("R9").Value = 0.1233480176211454
$ws.Range("S9").Value = 0.4317180616740088

$ws.Range("B10").Value = 0.119496855345912
$ws.Range("D10").Value = 0.0220125786163522
$ws.Range("F10").Value = 0.07232704402515723
$ws.Range("J10").Value = 0.1077044025157233
$ws.Range("O10").Value = 0.01336477987421384
$ws.Range("Q10").Value = 0.1957547169811321
$ws.Range("R10").Value = 0.07311320754716981
$ws.Range("S10").Value = 0.3962264150943396

$ws.Range("G11").Value = 0.1432835820895522
$ws.Range("J11").Value = 0.06567164179104477
$ws.Range("K11").Value = 0.1850746268656716
$ws.Range("L11").Value = 0.6
$ws.Range("S11").Value = 0.005970149253731343

$ws.Range("G12").Value = 0.7227722772277227
$ws.Range("J12").Value = 0.2079207920792079
$ws.Range("K12").Value = 0.01485148514851485
$ws.Range("L12").Value = 0.02475247524752475
$ws.Range("S12").Value = 0.0297029702970297

$ws.Range("G13").Value = 0.5833333333333334
$ws.Range("J13").Value = 0.3541666666666667
$ws.Range("S13").Value = 0.0625

$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333

$ws.Range("F15").Value = 0.03092783505154639
$ws.Range("H15").Value = 0.1443298969072165
$ws.Range("I15").Value = 0.1082474226804124
$ws.Range("J15").Value = 0.3659793814432989
$ws.Range("K15").Value = 0.03608247422680412
$ws.Range("M15").Value = 0.005154639175257732
$ws.Range("N15").Value = 0.005154639175257732
$ws.Range("O15").Value = 0.04639175257731959
$ws.Range("S15").Value = 0.2577319587628866

$ws.Range("F16").Value = 0.004878048780487805
$ws.Range("H16").Value = 0.08292682926829269
$ws.Range("I16").Value = 0.1073170731707317
$ws.Range("J16").Value = 0.4048780487804878
$ws.Range("K16").Value = 0.1414634146341463
$ws.Range("M16").Value = 0.03414634146341464
$ws.Range("O16").Value = 0.05853658536585366
$ws.Range("S16").Value = 0.1658536585365854

$ws.Range("F17").Value = 0.02457002457002457
$ws.Range("H17").Value = 0.1597051597051597
$ws.Range("I17").Value = 0.085995085995086
$ws.Range("J17").Value = 0.3832923832923833
$ws.Range("K17").Value = 0.1081081081081081
$ws.Range("M17").Value = 0.02457002457002457
$ws.Range("N17").Value = 0.002457002457002457
$ws.Range("O17").Value = 0.05405405405405406
$ws.Range("S17").Value = 0.1572481572481572

$ws.Range("F18").Value = 0.02173913043478261
$ws.Range("H18").Value = 0.1630434782608696
$ws.Range("I18").Value = 0.09782608695652174
$ws.Range("J18").Value = 0.4021739130434783
$ws.Range("K18").Value = 0.1032608695652174
$ws.Range("M18").Value = 0.01630434782608696
$ws.Range("O18").Value = 0.06521739130434782
$ws.Range("S18").Value = 0.1304347826086956

$ws.Range("F19").Value = 0.01538461538461539
$ws.Range("H19").Value = 0.2090909090909091
$ws.Range("I19").Value = 0.09300699300699301
$ws.Range("J19").Value = 0.3461538461538461
$ws.Range("K19").Value = 0.1174825174825175
$ws.Range("M19").Value = 0.01958041958041958
$ws.Range("N19").Value = 0.001398601398601399
$ws.Range("O19").Value = 0.05314685314685315
$ws.Range("S19").Value = 0.1447552447552448
